$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting rows 84:89 down to 85:90
$ws.Rows.Item(84).Insert()

# Fill the new row 84 with data (copy of old row84 values, with updated fields)
$ws.Cells.Item(84, 1).Value = 7
$ws.Cells.Item(84, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(84, 3).Value = "Ñuble"
$ws.Cells.Item(84, 4).Value = 44776
$ws.Cells.Item(84, 5).Value = 16
$ws.Cells.Item(84, 6).Value = 100112031
$ws.Cells.Item(84, 7).Value = "Poroto verde"
$ws.Cells.Item(84, 8).Value = "Magnum"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 30
$ws.Cells.Item(84, 11).Value = 30000
$ws.Cells.Item(84, 12).Value = 30000
$ws.Cells.Item(84, 13).Value = 30000
$ws.Cells.Item(84, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(84, 15).Value = "Perú"
$ws.Cells.Item(84, 16).Value = 1200
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"

# Apply the date style (numFmtId 165) to D84, matching the other date cells in column D
$ws.Cells.Item(84, 4).NumberFormat = $ws.Cells.Item(85, 4).NumberFormat
